$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 4.063713
$ws.Cells.Item(2, 8).Value = 12.191139
$ws.Cells.Item(2, 9).Value = 0.5065008440615062
$ws.Cells.Item(2, 10).Value = 0.5065008440615063
$ws.Cells.Item(2, 13).Value = 3.478915
$ws.Cells.Item(2, 14).Value = 10.436745
$ws.Cells.Item(2, 15).Value = 0.2270489597131295
$ws.Cells.Item(2, 16).Value = 0.2270489597131294
$ws.Cells.Item(2, 17).Value = 14.137312111395
$ws.Cells.Item(2, 18).Value = 127.235809002555
$ws.Cells.Item(2, 19).Value = 0.115000489737987
$ws.Cells.Item(2, 20).Value = 0.115000489737987
$ws.Cells.Item(3, 7).Value = 4.063713
$ws.Cells.Item(3, 8).Value = 12.191139
$ws.Cells.Item(3, 9).Value = 0.5065008440615062
$ws.Cells.Item(3, 10).Value = 0.5065008440615063
$ws.Cells.Item(3, 15).Value = 0.6268017733442142
$ws.Cells.Item(3, 16).Value = 0.6268017733442142
$ws.Cells.Item(3, 17).Value = 39.028112319647
$ws.Cells.Item(3, 18).Value = 351.253010876823
$ws.Cells.Item(3, 19).Value = 0.3174756272580934
$ws.Cells.Item(3, 20).Value = 0.3174756272580934
$ws.Cells.Item(4, 7).Value = 4.063713
$ws.Cells.Item(4, 8).Value = 12.191139
$ws.Cells.Item(4, 9).Value = 0.5065008440615062
$ws.Cells.Item(4, 10).Value = 0.5065008440615063
$ws.Cells.Item(4, 13).Value = 1.465966333333333
$ws.Cells.Item(4, 14).Value = 4.397899
$ws.Cells.Item(4, 15).Value = 0.09567526971995696
$ws.Cells.Item(4, 16).Value = 0.09567526971995696
$ws.Cells.Item(4, 17).Value = 5.957266446328999
$ws.Cells.Item(4, 18).Value = 53.615398016961
$ws.Cells.Item(4, 19).Value = 0.04845960486897047
$ws.Cells.Item(4, 20).Value = 0.04845960486897048
$ws.Cells.Item(5, 7).Value = 4.063713
$ws.Cells.Item(5, 8).Value = 12.191139
$ws.Cells.Item(5, 9).Value = 0.5065008440615062
$ws.Cells.Item(5, 10).Value = 0.5065008440615063
$ws.Cells.Item(5, 13).Value = 0.7733783333333334
$ws.Cells.Item(5, 14).Value = 2.320135
$ws.Cells.Item(5, 15).Value = 0.05047399722269938
$ws.Cells.Item(5, 16).Value = 0.05047399722269937
$ws.Cells.Item(5, 17).Value = 3.142787587085
$ws.Cells.Item(5, 18).Value = 28.285088283765
$ws.Cells.Item(5, 19).Value = 0.02556512219645535
$ws.Cells.Item(5, 20).Value = 0.02556512219645536
$ws.Cells.Item(6, 9).Value = 0.2604012840237886
$ws.Cells.Item(6, 10).Value = 0.2604012840237886
$ws.Cells.Item(6, 13).Value = 3.478915
$ws.Cells.Item(6, 14).Value = 10.436745
$ws.Cells.Item(6, 15).Value = 0.2270489597131295
$ws.Cells.Item(6, 16).Value = 0.2270489597131294
$ws.Cells.Item(6, 17).Value = 7.268248946896667
$ws.Cells.Item(6, 18).Value = 65.41424052207
$ws.Cells.Item(6, 19).Value = 0.05912384064556436
$ws.Cells.Item(6, 20).Value = 0.05912384064556436
$ws.Cells.Item(7, 9).Value = 0.2604012840237886
$ws.Cells.Item(7, 10).Value = 0.2604012840237886
$ws.Cells.Item(7, 15).Value = 0.6268017733442142
$ws.Cells.Item(7, 16).Value = 0.6268017733442142
$ws.Cells.Item(7, 19).Value = 0.1632199866072211
$ws.Cells.Item(7, 20).Value = 0.1632199866072211
$ws.Cells.Item(8, 9).Value = 0.2604012840237886
$ws.Cells.Item(8, 10).Value = 0.2604012840237886
$ws.Cells.Item(8, 13).Value = 1.465966333333333
$ws.Cells.Item(8, 14).Value = 4.397899
$ws.Cells.Item(8, 15).Value = 0.09567526971995696
$ws.Cells.Item(8, 16).Value = 0.09567526971995696
$ws.Cells.Item(8, 17).Value = 3.062738887968222
$ws.Cells.Item(8, 18).Value = 27.564649991714
$ws.Cells.Item(8, 19).Value = 0.02491396308439909
$ws.Cells.Item(8, 20).Value = 0.0249139630843991
$ws.Cells.Item(9, 9).Value = 0.2604012840237886
$ws.Cells.Item(9, 10).Value = 0.2604012840237886
$ws.Cells.Item(9, 13).Value = 0.7733783333333334
$ws.Cells.Item(9, 14).Value = 2.320135
$ws.Cells.Item(9, 15).Value = 0.05047399722269938
$ws.Cells.Item(9, 16).Value = 0.05047399722269937
$ws.Cells.Item(9, 17).Value = 1.615764184178889
$ws.Cells.Item(9, 18).Value = 14.54187765761
$ws.Cells.Item(9, 19).Value = 0.01314349368660406
$ws.Cells.Item(9, 20).Value = 0.01314349368660406
$ws.Cells.Item(10, 7).Value = 1.588356333333333
$ws.Cells.Item(10, 8).Value = 4.765069
$ws.Cells.Item(10, 9).Value = 0.1979725988286506
$ws.Cells.Item(10, 10).Value = 0.1979725988286507
$ws.Cells.Item(10, 13).Value = 3.478915
$ws.Cells.Item(10, 14).Value = 10.436745
$ws.Cells.Item(10, 15).Value = 0.2270489597131295
$ws.Cells.Item(10, 16).Value = 0.2270489597131294
$ws.Cells.Item(10, 17).Value = 5.525756673378333
$ws.Cells.Item(10, 18).Value = 49.731810060405
$ws.Cells.Item(10, 19).Value = 0.04494947261574984
$ws.Cells.Item(10, 20).Value = 0.04494947261574984
$ws.Cells.Item(11, 7).Value = 1.588356333333333
$ws.Cells.Item(11, 8).Value = 4.765069
$ws.Cells.Item(11, 9).Value = 0.1979725988286506
$ws.Cells.Item(11, 10).Value = 0.1979725988286507
$ws.Cells.Item(11, 15).Value = 0.6268017733442142
$ws.Cells.Item(11, 16).Value = 0.6268017733442142
$ws.Cells.Item(11, 17).Value = 15.25465734931478
$ws.Cells.Item(11, 18).Value = 137.291916143833
$ws.Cells.Item(11, 19).Value = 0.1240895760193609
$ws.Cells.Item(11, 20).Value = 0.1240895760193609
$ws.Cells.Item(12, 7).Value = 1.588356333333333
$ws.Cells.Item(12, 8).Value = 4.765069
$ws.Cells.Item(12, 9).Value = 0.1979725988286506
$ws.Cells.Item(12, 10).Value = 0.1979725988286507
$ws.Cells.Item(12, 13).Value = 1.465966333333333
$ws.Cells.Item(12, 14).Value = 4.397899
$ws.Cells.Item(12, 15).Value = 0.09567526971995696
$ws.Cells.Item(12, 16).Value = 0.09567526971995696
$ws.Cells.Item(12, 17).Value = 2.328476910003444
$ws.Cells.Item(12, 18).Value = 20.956292190031
$ws.Cells.Item(12, 19).Value = 0.01894108179009199
$ws.Cells.Item(12, 20).Value = 0.01894108179009199
$ws.Cells.Item(13, 7).Value = 1.588356333333333
$ws.Cells.Item(13, 8).Value = 4.765069
$ws.Cells.Item(13, 9).Value = 0.1979725988286506
$ws.Cells.Item(13, 10).Value = 0.1979725988286507
$ws.Cells.Item(13, 13).Value = 0.7733783333333334
$ws.Cells.Item(13, 14).Value = 2.320135
$ws.Cells.Item(13, 15).Value = 0.05047399722269938
$ws.Cells.Item(13, 16).Value = 0.05047399722269937
$ws.Cells.Item(13, 17).Value = 1.228400373812778
$ws.Cells.Item(13, 18).Value = 11.055603364315
$ws.Cells.Item(13, 19).Value = 0.009992468403447891
$ws.Cells.Item(13, 20).Value = 0.009992468403447891
$ws.Cells.Item(14, 7).Value = 0.281814
$ws.Cells.Item(14, 8).Value = 0.845442
$ws.Cells.Item(14, 9).Value = 0.03512527308605438
$ws.Cells.Item(14, 10).Value = 0.03512527308605439
$ws.Cells.Item(14, 13).Value = 3.478915
$ws.Cells.Item(14, 14).Value = 10.436745
$ws.Cells.Item(14, 15).Value = 0.2270489597131295
$ws.Cells.Item(14, 16).Value = 0.2270489597131294
$ws.Cells.Item(14, 17).Value = 0.98040695181
$ws.Cells.Item(14, 18).Value = 8.82366256629
$ws.Cells.Item(14, 19).Value = 0.007975156713828233
$ws.Cells.Item(14, 20).Value = 0.007975156713828233
$ws.Cells.Item(15, 7).Value = 0.281814
$ws.Cells.Item(15, 8).Value = 0.845442
$ws.Cells.Item(15, 9).Value = 0.03512527308605438
$ws.Cells.Item(15, 10).Value = 0.03512527308605439
$ws.Cells.Item(15, 15).Value = 0.6268017733442142
$ws.Cells.Item(15, 16).Value = 0.6268017733442142
$ws.Cells.Item(15, 17).Value = 2.706556404266
$ws.Cells.Item(15, 18).Value = 24.359007638394
$ws.Cells.Item(15, 19).Value = 0.02201658345953869
$ws.Cells.Item(15, 20).Value = 0.02201658345953869
$ws.Cells.Item(16, 7).Value = 0.281814
$ws.Cells.Item(16, 8).Value = 0.845442
$ws.Cells.Item(16, 9).Value = 0.03512527308605438
$ws.Cells.Item(16, 10).Value = 0.03512527308605439
$ws.Cells.Item(16, 13).Value = 1.465966333333333
$ws.Cells.Item(16, 14).Value = 4.397899
$ws.Cells.Item(16, 15).Value = 0.09567526971995696
$ws.Cells.Item(16, 16).Value = 0.09567526971995696
$ws.Cells.Item(16, 17).Value = 0.413129836262
$ws.Cells.Item(16, 18).Value = 3.718168526358
$ws.Cells.Item(16, 19).Value = 0.003360619976495398
$ws.Cells.Item(16, 20).Value = 0.003360619976495399
$ws.Cells.Item(17, 7).Value = 0.281814
$ws.Cells.Item(17, 8).Value = 0.845442
$ws.Cells.Item(17, 9).Value = 0.03512527308605438
$ws.Cells.Item(17, 10).Value = 0.03512527308605439
$ws.Cells.Item(17, 13).Value = 0.7733783333333334
$ws.Cells.Item(17, 14).Value = 2.320135
$ws.Cells.Item(17, 15).Value = 0.05047399722269938
$ws.Cells.Item(17, 16).Value = 0.05047399722269937
$ws.Cells.Item(17, 17).Value = 0.21794884163
$ws.Cells.Item(17, 18).Value = 1.96153957467
$ws.Cells.Item(17, 19).Value = 0.001772912936192066
$ws.Cells.Item(17, 20).Value = 0.001772912936192066
